# Generate Report for Handback
# - Flip Status cells from "Ready for handoff" to "Handed back: in sync with en-US"
# - Record the handback target/handback-file columns (I/J) with a hyperlink to the
#   source doc and the generated .xlf name, and stamp the handback datetime (K)
# - Widen a few columns that now hold longer text

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$statusText = "Handed back: in sync with en-US"
$hyperlinkUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/15b7e3bef064435f7bb491c716c88126b4ff8845/e2e/28d88ecf-cec1-483d-8822-62f010084577.md"
$hyperlinkDisplay = "28d88ecf-cec1-483d-8822-62f010084577.md"
$hyperlinkColor = 15570276

# ---------------------------------------------------------------------------
# Overview sheet: Status columns for zh-cn (E) and de-de (F)
# ---------------------------------------------------------------------------
$overview.Range("E2").Value = $statusText
$overview.Range("F2").Value = $statusText
$overview.Range("E3").Value = $statusText
$overview.Range("F3").Value = $statusText

$overview.Columns.Item(5).ColumnWidth = 29.166666666666664
$overview.Columns.Item(6).ColumnWidth = 29.166666666666664

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$zhcn.Range("C2").Value = $statusText
$zhcn.Range("C3").Value = $statusText

$zhcn.Range("I2").Value = $hyperlinkDisplay
$zhcn.Range("I3").Value = $hyperlinkDisplay
$zhcn.Range("J2").Value = "28d88ecf-cec1-483d-8822-62f010084577.80c718f14bd74f96b3e433d7807a7ee00e2a0328.zh-cn.xlf"
$zhcn.Range("J3").Value = "28d88ecf-cec1-483d-8822-62f010084577.80c718f14bd74f96b3e433d7807a7ee00e2a0328.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-09-04 03:07:51"
$zhcn.Range("K3").Value = "2016-09-04 03:07:51"

$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $hyperlinkUrl, "", "", $hyperlinkDisplay)
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), $hyperlinkUrl, "", "", $hyperlinkDisplay)

$zhcn.Range("I2").Font.Underline = $true
$zhcn.Range("I2").Font.Color = $hyperlinkColor
$zhcn.Range("I3").Font.Underline = $true
$zhcn.Range("I3").Font.Color = $hyperlinkColor

$zhcn.Columns.Item(3).ColumnWidth = 29.166666666666664
$zhcn.Columns.Item(9).ColumnWidth = 39.166666666666664
$zhcn.Columns.Item(10).ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$dede.Range("C2").Value = $statusText
$dede.Range("C3").Value = $statusText

$dede.Range("I2").Value = $hyperlinkDisplay
$dede.Range("I3").Value = $hyperlinkDisplay
$dede.Range("J2").Value = "28d88ecf-cec1-483d-8822-62f010084577.80c718f14bd74f96b3e433d7807a7ee00e2a0328.de-de.xlf"
$dede.Range("J3").Value = "28d88ecf-cec1-483d-8822-62f010084577.80c718f14bd74f96b3e433d7807a7ee00e2a0328.de-de.xlf"
$dede.Range("K2").Value = "2016-09-04 03:07:59"
$dede.Range("K3").Value = "2016-09-04 03:07:59"

$dede.Hyperlinks.Add($dede.Range("I2"), $hyperlinkUrl, "", "", $hyperlinkDisplay)
$dede.Hyperlinks.Add($dede.Range("I3"), $hyperlinkUrl, "", "", $hyperlinkDisplay)

$dede.Range("I2").Font.Underline = $true
$dede.Range("I2").Font.Color = $hyperlinkColor
$dede.Range("I3").Font.Underline = $true
$dede.Range("I3").Font.Color = $hyperlinkColor

$dede.Columns.Item(3).ColumnWidth = 29.166666666666664
$dede.Columns.Item(9).ColumnWidth = 39.166666666666664
$dede.Columns.Item(10).ColumnWidth = 39.166666666666664
